# Auto-generated edit script applying the diff to Zalera_Profits workbook
# Updates currentAveragePrice / Leve price / profit columns (H:N) on several
# rows across the ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR sheets, matching a
# scheduled market-data refresh.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 103.73684
$ws.Range("I28").Value = 107
$ws.Range("K28").Value = 107
$ws.Range("M28").Value = 378

$ws.Range("H33").Value = 298.73685
$ws.Range("I33").Value = 304.8125
$ws.Range("J33").Value = 266.33334
$ws.Range("K33").Value = 304.8125
$ws.Range("L33").Value = 266.33334
$ws.Range("M33").Value = -75.8125
$ws.Range("N33").Value = -724.33334

$ws.Range("H86").Value = 12102.167
$ws.Range("I86").Value = 13350.223
$ws.Range("J86").Value = 8358
$ws.Range("K86").Value = 13350.223
$ws.Range("L86").Value = 8358
$ws.Range("M86").Value = -12227.223
$ws.Range("N86").Value = -10604

$ws.Range("H89").Value = 12102.167
$ws.Range("I89").Value = 13350.223
$ws.Range("J89").Value = 8358
$ws.Range("K89").Value = 66751.11500000001
$ws.Range("L89").Value = 41790
$ws.Range("M89").Value = -61135.11500000001
$ws.Range("N89").Value = -53022

$ws.Range("H123").Value = 119993
$ws.Range("J123").Value = 119993
$ws.Range("L123").Value = 119993
$ws.Range("N123").Value = -129793

$ws.Range("H141").Value = 2027.9584
$ws.Range("I141").Value = 1757.8636
$ws.Range("K141").Value = 5273.5908
$ws.Range("M141").Value = -93.59079999999994

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1565.8125
$ws.Range("I45").Value = 1523.5454
$ws.Range("K45").Value = 1523.5454
$ws.Range("M45").Value = -1146.5454

$ws.Range("H103").Value = 38340.5
$ws.Range("J103").Value = 38340.5
$ws.Range("L103").Value = 38340.5
$ws.Range("N103").Value = -40684.5

$ws.Range("H132").Value = 13145.0625
$ws.Range("I132").Value = 4077.45
$ws.Range("K132").Value = 12232.35
$ws.Range("M132").Value = -9702.349999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 2447.0476
$ws.Range("I107").Value = 1915.2632
$ws.Range("K107").Value = 1915.2632
$ws.Range("M107").Value = 4.736799999999903

$ws.Range("H132").Value = 111917.375
$ws.Range("J132").Value = 111917.375
$ws.Range("L132").Value = 111917.375
$ws.Range("N132").Value = -122037.375

$ws.Range("H134").Value = 4530.9062
$ws.Range("I134").Value = 3580.0715
$ws.Range("J134").Value = 11186.75
$ws.Range("K134").Value = 10740.2145
$ws.Range("L134").Value = 33560.25
$ws.Range("M134").Value = -8205.2145
$ws.Range("N134").Value = -38630.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 1600.4286
$ws.Range("J2").Value = 2000
$ws.Range("L2").Value = 2000
$ws.Range("N2").Value = -2226

$ws.Range("H31").Value = 7758.391
$ws.Range("I31").Value = 5768.8184
$ws.Range("K31").Value = 5768.8184
$ws.Range("M31").Value = -5473.8184

$ws.Range("H34").Value = 7758.391
$ws.Range("I34").Value = 5768.8184
$ws.Range("K34").Value = 5768.8184
$ws.Range("M34").Value = -5566.8184

$ws.Range("H58").Value = 5397.5386
$ws.Range("I58").Value = 2463.1428
$ws.Range("J58").Value = 8821
$ws.Range("K58").Value = 2463.1428
$ws.Range("L58").Value = 8821
$ws.Range("M58").Value = -2260.1428
$ws.Range("N58").Value = -9227

$ws.Range("H94").Value = 13007.667
$ws.Range("I94").Value = 16756.666
$ws.Range("K94").Value = 16756.666
$ws.Range("M94").Value = -16305.666

$ws.Range("H107").Value = 320.25
$ws.Range("I107").Value = 325.36365
$ws.Range("J107").Value = 264
$ws.Range("K107").Value = 325.36365
$ws.Range("L107").Value = 264
$ws.Range("M107").Value = 1594.63635
$ws.Range("N107").Value = -4104

$ws.Range("H134").Value = 2381.8235
$ws.Range("J134").Value = 7285
$ws.Range("L134").Value = 21855
$ws.Range("N134").Value = -26925

$ws.Range("H136").Value = 5397.5386
$ws.Range("I136").Value = 2463.1428
$ws.Range("J136").Value = 8821
$ws.Range("K136").Value = 7389.428400000001
$ws.Range("L136").Value = 26463
$ws.Range("M136").Value = -4839.428400000001
$ws.Range("N136").Value = -31563

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H117").Value = 1420.4
$ws.Range("I117").Value = 386.75
$ws.Range("J117").Value = 5555
$ws.Range("K117").Value = 1160.25
$ws.Range("L117").Value = 16665
$ws.Range("M117").Value = 2281.75
$ws.Range("N117").Value = -23549

$ws.Range("H129").Value = 22729094
$ws.Range("I129").Value = 690.25
$ws.Range("J129").Value = 50003176
$ws.Range("K129").Value = 2070.75
$ws.Range("L129").Value = 150009528
$ws.Range("M129").Value = 2929.25
$ws.Range("N129").Value = -150019528

$ws.Range("H131").Value = 30327888
$ws.Range("I131").Value = 166667660
$ws.Range("K131").Value = 500002980
$ws.Range("M131").Value = -499997940

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H23").Value = 1007.5
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 1007.5
$ws.Range("K23").Value = 0
$ws.Range("L23").ClearContents()
$ws.Range("M23").Value = 1007.5
$ws.Range("N23").Value = -1453.5

$ws.Range("H46").Value = 49083.09
$ws.Range("J46").Value = 56101.555
$ws.Range("L46").Value = 56101.555
$ws.Range("N46").Value = -56413.555

$ws.Range("H57").Value = 59748.5
$ws.Range("I57").Value = 56333
$ws.Range("K57").Value = 56333
$ws.Range("M57").Value = -55513

$ws.Range("H70").Value = 6998.6665
$ws.Range("I70").Value = 6998.6665
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 6998.6665
$ws.Range("L70").ClearContents()
$ws.Range("M70").Value = -6728.6665
$ws.Range("N70").Value = 0

$ws.Range("H73").Value = 6998.6665
$ws.Range("I73").Value = 6998.6665
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 6998.6665
$ws.Range("L73").ClearContents()
$ws.Range("M73").Value = -6062.6665
$ws.Range("N73").Value = 0

$ws.Range("H93").Value = 45039.168
$ws.Range("J93").Value = 45039.168
$ws.Range("L93").Value = 45039.168
$ws.Range("N93").Value = -48783.168

$ws.Range("H107").Value = 3069.5
$ws.Range("I107").Value = 759.3333
$ws.Range("K107").Value = 759.3333
$ws.Range("M107").Value = 1160.6667

$ws.Range("H113").Value = 2724.5
$ws.Range("I113").Value = 2724.5
$ws.Range("K113").Value = 2724.5
$ws.Range("M113").Value = -554.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1062.0857
$ws.Range("I16").Value = 1051.0968
$ws.Range("K16").Value = 1051.0968
$ws.Range("M16").Value = -881.0968

$ws.Range("H46").Value = 3069.65
$ws.Range("I46").Value = 1568.6154
$ws.Range("J46").Value = 5857.2856
$ws.Range("K46").Value = 1568.6154
$ws.Range("L46").Value = 5857.2856
$ws.Range("M46").Value = -1380.6154
$ws.Range("N46").Value = -6233.2856

$ws.Range("H104").Value = 23999.5
$ws.Range("J104").Value = 23999.5
$ws.Range("L104").Value = 23999.5
$ws.Range("N104").Value = -30987.5

$ws.Range("H122").Value = 20837990
$ws.Range("I122").Value = 27782550
$ws.Range("K122").Value = 83347650
$ws.Range("M122").Value = -83345200

$ws.Range("H132").Value = 5735.048
$ws.Range("I132").Value = 3702.5334
$ws.Range("K132").Value = 11107.6002
$ws.Range("M132").Value = -8577.600199999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 3783.0667
$ws.Range("I122").Value = 3764.75
$ws.Range("K122").Value = 11294.25
$ws.Range("M122").Value = -8844.25

$ws.Range("H123").Value = 49999.5
$ws.Range("J123").Value = 49999.5
$ws.Range("L123").Value = 49999.5
$ws.Range("N123").Value = -59799.5

$ws.Range("H126").Value = 251726.25
$ws.Range("I126").Value = 1000000
$ws.Range("J126").Value = 2301.6667
$ws.Range("K126").Value = 3000000
$ws.Range("L126").Value = 6905.000100000001
$ws.Range("M126").Value = -2997530
$ws.Range("N126").Value = -11845.0001

$ws.Range("H132").Value = 5793.3184
$ws.Range("I132").Value = 4643.7036
$ws.Range("J132").Value = 7619.1763
$ws.Range("K132").Value = 13931.1108
$ws.Range("L132").Value = 22857.5289
$ws.Range("M132").Value = -11401.1108
$ws.Range("N132").Value = -27917.5289
